$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.590432228283676
$ws.Range("C2").Value = 0.601344860710855
$ws.Range("D2").Value = 0.613245033112583
$ws.Range("E2").Value = 0.587147030185005
$ws.Range("F2").Value = 0.640070921985816
$ws.Range("G2").Value = 0.58043758043758
